$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.424.03'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '2.371.72'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.32'
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.46'
$ws.Range("E6").Value = '  -2.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.632'
$ws.Range("E7").Value = '  -0.83%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.614'
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.93'
$ws.Range("E10").Value = '  -2.97%  '
$ws.Range("E11").Value = '  -1.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.51'
$ws.Range("E12").Value = '  -1.89%  '
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.985'
$ws.Range("E14").Value = '  -3.19%  '
$ws.Range("D15").Value = '2.731.78'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.38'
$ws.Range("E16").Value = '  -2.62%  '
$ws.Range("D17").Value = '2.369.96'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").Value = '45.395.41'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.59'
$ws.Range("E19").Value = '  +18.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.29'
$ws.Range("E20").Value = '  -5.07%  '
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.42'
$ws.Range("E22").Value = '  -2.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.59'
$ws.Range("E23").Value = '  +1.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '261.48'
$ws.Range("E24").Value = '  -3.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.15'
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.46'
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.30'
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0967'
$ws.Range("E30").Value = '  +1.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.34'
$ws.Range("E31").Value = '  -2.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.17'
$ws.Range("E32").Value = '  -3.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '166.67'
$ws.Range("E33").Value = '  -1.92%  '
$ws.Range("E34").Value = '  -3.10%  '
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.70'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.90'
$ws.Range("E38").Value = '  +9.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.95'
$ws.Range("E39").Value = '  +0.48%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.94'
$ws.Range("E40").Value = '  -4.30%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0355'
$ws.Range("E41").Value = '  -2.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.81'
$ws.Range("E42").Value = '  -6.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.48'
$ws.Range("E43").Value = '  -1.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.05'
$ws.Range("E44").Value = '  -3.54%  '
$ws.Range("E45").Value = '  -5.67%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.95'
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.811.80'
$ws.Range("E48").Value = '  +9.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.09'
$ws.Range("E49").Value = '  +5.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.19'
$ws.Range("E50").Value = '  -6.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.21'
$ws.Range("E51").Value = '  -0.94%  '
